# post-stats analysis and thesis drafting
#
# Removes the standalone "Use R to:" lead-in paragraph from the
# "Evaluate Satellite Data Accuracy" callout box, and tightens the
# wording of the first bullet from
# "Load in the comprehensive statistical dataset." to
# "Load comprehensive statistical dataset."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Rounded Rectangle 389" shape that holds the callout text.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "Rounded Rectangle 389") {
        $shape = $s.Shapes.Item($i)
        break
    }
}

$tr = $shape.TextFrame.TextRange

# Paragraph layout before the edit:
#   1. "Evaluate Satellite Data Accuracy"
#   2. "Use R to:"
#   3. "Load in the comprehensive statistical dataset."
#   4. "Generate plots that compare ..."
#
# Drop paragraph 2 entirely.
$tr.Paragraphs(2, 1).Delete()

# Paragraph 3 is now paragraph 2. Reword it while keeping its bullet
# formatting and run properties intact. Setting the text directly would
# keep the common "Load " prefix as a separate run, so stage the change
# through an unrelated placeholder first to force a clean single run.
$bulletPara = $tr.Paragraphs(2, 1)
$bulletPara.Text = "zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$tr.Paragraphs(2, 1).Text = "Load comprehensive statistical dataset."
